$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.962.82'
$ws.Range('E2').Value = '  -4.39%  '
$ws.Range('D3').Value = '2.242.96'
$ws.Range('E3').Value = '  -4.69%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '232.22'
$ws.Range('E5').Value = '  -3.71%  '
$ws.Range('D6').Value = '0.631'
$ws.Range('E6').Value = '  -6.47%  '
$ws.Range('D7').Value = '69.32'
$ws.Range('E7').Value = '  -5.45%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').Value = '0.553'
$ws.Range('E9').Value = '  -8.59%  '
$ws.Range('D10').Value = '0.0983'
$ws.Range('E10').Value = '  -2.33%  '
$ws.Range('D11').Value = '57.86'
$ws.Range('E11').Value = '  -2.50%  '
$ws.Range('D12').Value = '35.51'
$ws.Range('E12').Value = '  +5.84%  '
$ws.Range('E13').Value = '  -3.50%  '
$ws.Range('D14').Value = '6.74'
$ws.Range('E14').Value = '  -7.90%  '
$ws.Range('D15').Value = '2.577.21'
$ws.Range('E15').Value = '  -4.74%  '
$ws.Range('D16').Value = '14.97'
$ws.Range('E16').Value = '  -8.85%  '
$ws.Range('D17').Value = '0.857'
$ws.Range('E17').Value = '  -5.34%  '
$ws.Range('D18').Value = '2.244.25'
$ws.Range('E18').Value = '  -4.79%  '
$ws.Range('D19').Value = '41.854.87'
$ws.Range('E19').Value = '  -4.44%  '
$ws.Range('D20').Value = '0.0₃0967'
$ws.Range('E20').Value = '  -6.28%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').Value = '73.17'
$ws.Range('E21').Value = '  -5.65%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '6.22'
$ws.Range('E22').Value = '  -7.40%  '
$ws.Range('D23').Value = '236.04'
$ws.Range('E23').Value = '  -7.93%  '
$ws.Range('D24').Value = '2.04'
$ws.Range('E24').Value = '  +5.21%  '
$ws.Range('E25').Value = '  -0.10%  '
$ws.Range('D26').Value = '3.64'
$ws.Range('E26').Value = '  -2.34%  '
$ws.Range('D27').Value = '2.35'
$ws.Range('E27').Value = '  -6.16%  '
$ws.Range('D28').Value = '9.99'
$ws.Range('E28').Value = '  -5.66%  '
$ws.Range('E29').Value = '  -1.91%  '
$ws.Range('D30').Value = '168.73'
$ws.Range('E30').Value = '  -4.83%  '
$ws.Range('D31').Value = '20.55'
$ws.Range('E31').Value = '  -9.08%  '
$ws.Range('D32').Value = '0.118'
$ws.Range('E32').Value = '  -7.69%  '
$ws.Range('D33').Value = '0.127'
$ws.Range('E33').Value = '  -7.42%  '
$ws.Range('D34').Value = '5.46'
$ws.Range('E34').Value = '  -0.29%  '
$ws.Range('D35').Value = '0.0712'
$ws.Range('E35').Value = '  -5.73%  '
$ws.Range('D36').Value = '4.76'
$ws.Range('E36').Value = '  -7.96%  '
$ws.Range('D37').Value = '3.59'
$ws.Range('E37').Value = '  -5.45%  '
$ws.Range('D38').Value = '22.11'
$ws.Range('E38').Value = '  +16.53%  '
$ws.Range('D39').Value = '2.25'
$ws.Range('E39').Value = '  -5.70%  '
$ws.Range('D40').Value = '6.02'
$ws.Range('E40').Value = '  -6.54%  '
$ws.Range('D41').Value = '0.0265'
$ws.Range('E41').Value = '  -4.22%  '
$ws.Range('D42').Value = '66.59'
$ws.Range('E42').Value = '  -1.46%  '
$ws.Range('D43').Value = '5.00'
$ws.Range('E43').Value = '  -2.52%  '
$ws.Range('D44').Value = '9.07'
$ws.Range('E44').Value = '  -1.94%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = '0.100'
$ws.Range('E45').Value = '  -8.42%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = '0.190'
$ws.Range('E46').Value = '  -6.00%  '
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('B48').Value = 'SynthetixNetwork'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D48').Value = '4.33'
$ws.Range('E48').Value = '  +6.84%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '2.34'
$ws.Range('E49').Value = '  -6.85%  '
$ws.Range('D50').Value = '1.17'
$ws.Range('E50').Value = '  -7.11%  '
$ws.Range('D51').Value = '9.87'
$ws.Range('E51').Value = '  +3.03%  '
